$d = $word.ActiveDocument

# Change 1: "2 întrebări ușoare de teorie de la cursuri" -> "2 exemple de cod de la cursuri"
$d.Content.Find.Execute("2 întrebări ușoare de teorie de la cursuri", $true, $false, $false, $false, $false, $true, 1, $false, "2 exemple de cod de la cursuri", 2)

# Change 2: insert " sau laboratoare" after the "(mai puțin primul și ultimele trei)" before ", la întâmplare;"
$d.Content.Find.Execute("trei), la întâmplare;", $true, $false, $false, $false, $false, $true, 1, $false, "trei) sau laboratoare, la întâmplare;", 2)

# Change 3: "2 întrebări mai grele de teorie de la cursuri" -> "2 întrebări de teorie de la cursuri"
$d.Content.Find.Execute("2 întrebări mai grele de teorie de la cursuri", $true, $false, $false, $false, $false, $true, 1, $false, "2 întrebări de teorie de la cursuri", 2)

# Change 4: "Dacă răspundeți bine la prima întrebare ușoară, " -> "Dacă explicați bine la primul exemplu de cod, "
$d.Content.Find.Execute("Dacă răspundeți bine la prima întrebare ușoară, ", $true, $false, $false, $false, $false, $true, 1, $false, "Dacă explicați bine la primul exemplu de cod, ", 2)

# Change 5: " puncte, și nu mai primiți altă întrebare ușoară de răspuns;" -> " puncte, și nu mai primiți alt exemplu de cod de explicat;"
$d.Content.Find.Execute(" puncte, și nu mai primiți altă întrebare ușoară de răspuns;", $true, $false, $false, $false, $false, $true, 1, $false, " puncte, și nu mai primiți alt exemplu de cod de explicat;", 2)

# Change 6: "Dacă nu, mai aveți o șansă, dar dacă răspundeți bine de data aceasta, " -> "Dacă nu, mai aveți o șansă, dar dacă explicați bine de data aceasta, "
$d.Content.Find.Execute("Dacă nu, mai aveți o șansă, dar dacă răspundeți bine de data aceasta, ", $true, $false, $false, $false, $false, $true, 1, $false, "Dacă nu, mai aveți o șansă, dar dacă explicați bine de data aceasta, ", 2)

# Change 7: "Dacă răspundeți bine la prima întrebare grea, " -> "Dacă răspundeți bine la prima întrebare de teorie, "
$d.Content.Find.Execute("Dacă răspundeți bine la prima întrebare grea, ", $true, $false, $false, $false, $false, $true, 1, $false, "Dacă răspundeți bine la prima întrebare de teorie, ", 2)

# Change 8: "puncte, și nu mai primiți altă întrebare grea de răspuns;" -> "puncte, și nu mai primiți altă întrebare de teorie de răspuns;"
$d.Content.Find.Execute("puncte, și nu mai primiți altă întrebare grea de răspuns;", $true, $false, $false, $false, $false, $true, 1, $false, "puncte, și nu mai primiți altă întrebare de teorie de răspuns;", 2)
